# fish_pond.xlsx - "debugging; next step is to debug lure on special branch."
#
# PondList sheet: two lake entries get renamed, and their "pond image" cell
# (column C) is repointed to reuse the lake's own name instead of a separate
# image-name string (RO_Lake/ogle_lake were using stale ro01_2D/ogle01_2D
# image ids in column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PondList")

# Row 5: RO_Lake -> Noob_Lake
$ws.Range("B5").Value = "Noob_Lake"
$ws.Range("B5").Copy()
$ws.Range("C5").PasteSpecial(-4122)          # xlPasteFormats: C5 drops the
                                              # old "image" fill/alignment and
                                              # picks up B5's plain style
$ws.Range("C5").Value = "Noob_Lake"

# Row 6: ogle_lake -> Ogle_Lake
$ws.Range("B6").Value = "Ogle_Lake"
$ws.Range("B6").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("C6").Value = "Ogle_Lake"

$excel.CutCopyMode = 0

# Re-point the "no duplicate values" validation on id/name columns A:B to a
# single contiguous range (it was split into three pieces before).
$rng = $ws.Range("A1:B1048576")
$rng.Validation.Delete()
$rng.Validation.Add(7, 2, 1, 'COUNTIF($A:$A,A1)<2')
$rng.Validation.ErrorTitle = "拒绝重复输入"
$rng.Validation.ErrorMessage = "当前输入的内容，与本区域的其他单元格内容重复。"
$rng.Validation.ShowInput = $false

# Leave the sheet with A5:C6 (the two rows just edited) selected.
$ws.Activate() | Out-Null
$ws.Range("A5:C6").Select() | Out-Null
